# Commit message: "handling empty cells at the end of rows"
#
# On sheet2, row 3 was A3=4, B3=5, C3=6. C3 is a trailing cell that should
# instead be empty/missing, so that the row exercises handling of rows whose
# final column(s) are empty. After this edit row 3 ends at B3.
$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# Drop the trailing/empty cell at the end of the row.
$ws2.Range("C3").ClearContents()

# The sheet's column formatting was only ever the sheet-wide default width,
# so it carries no real information; clearing cell formatting collapses that
# now-redundant column-width table.
$ws2.Cells.ClearFormats()
